# Add a new worksheet "ODI Batting Extra" as the last tab
# (mirrors: <sheet name="ODI Batting Extra" sheetId="4" state="visible" r:id="rId4"/>)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$headerSource = $wb.Worksheets.Item(1)      # "Player Info" - bold header style to copy
$codeSource = $wb.Worksheets.Item(2)        # "ODI Batting" - has MATCH_CODE "4619" as text

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row (row 1)
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Reuse the exact header style already used on the other sheets (bold font
# + border, centered/top aligned) by copying the format from an existing
# header cell instead of re-building it by hand.
$headerSource.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Data row (row 2) - MATCH_CODE for this match, stored as text like every
# other sheet in this workbook (copy the value from the existing "4619"
# text cell so it keeps its text type without picking up new formatting).
$codeSource.Range("D2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("F2").Value = "NO"

# Restore the originally active tab ("Player Info", index 0)
[void]$wb.Worksheets.Item(1).Select()
